$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 958.2
$ws.Cells.Item(18, 9).Value = 947.75
$ws.Cells.Item(18, 10).Value = 1000
$ws.Cells.Item(18, 11).Value = 947.75
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = -663.75
$ws.Cells.Item(18, 14).Value = -1568
$ws.Cells.Item(19, 8).Value = 169.72223
$ws.Cells.Item(19, 9).Value = 163.44444
$ws.Cells.Item(19, 10).Value = 176
$ws.Cells.Item(19, 11).Value = 163.44444
$ws.Cells.Item(19, 12).Value = 176
$ws.Cells.Item(19, 13).Value = 11.55556000000001
$ws.Cells.Item(19, 14).Value = -526
$ws.Cells.Item(40, 8).Value = 1216.8334
$ws.Cells.Item(40, 9).Value = 1060
$ws.Cells.Item(40, 11).Value = 1060
$ws.Cells.Item(40, 13).Value = -885
$ws.Cells.Item(62, 8).Value = 6433.3335
$ws.Cells.Item(62, 9).Value = 20005
$ws.Cells.Item(62, 10).Value = 3719
$ws.Cells.Item(62, 11).Value = 20005
$ws.Cells.Item(62, 12).Value = 3719
$ws.Cells.Item(62, 13).Value = -19381
$ws.Cells.Item(62, 14).Value = -4967
$ws.Cells.Item(65, 8).Value = 6433.3335
$ws.Cells.Item(65, 9).Value = 20005
$ws.Cells.Item(65, 10).Value = 3719
$ws.Cells.Item(65, 11).Value = 100025
$ws.Cells.Item(65, 12).Value = 18595
$ws.Cells.Item(65, 13).Value = -96905
$ws.Cells.Item(65, 14).Value = -24835
$ws.Cells.Item(74, 8).Value = 3409.9092
$ws.Cells.Item(74, 9).Value = 3140.8572
$ws.Cells.Item(74, 10).Value = 3880.75
$ws.Cells.Item(74, 11).Value = 3140.8572
$ws.Cells.Item(74, 12).Value = 3880.75
$ws.Cells.Item(74, 13).Value = -2204.8572
$ws.Cells.Item(74, 14).Value = -5752.75
$ws.Cells.Item(77, 8).Value = 3409.9092
$ws.Cells.Item(77, 9).Value = 3140.8572
$ws.Cells.Item(77, 10).Value = 3880.75
$ws.Cells.Item(77, 11).Value = 15704.286
$ws.Cells.Item(77, 12).Value = 19403.75
$ws.Cells.Item(77, 13).Value = -11024.286
$ws.Cells.Item(77, 14).Value = -28763.75
$ws.Cells.Item(94, 8).Value = 2524.1428
$ws.Cells.Item(94, 9).Value = 2112.9167
$ws.Cells.Item(94, 10).Value = 4991.5
$ws.Cells.Item(94, 11).Value = 2112.9167
$ws.Cells.Item(94, 12).Value = 4991.5
$ws.Cells.Item(94, 13).Value = -1661.9167
$ws.Cells.Item(94, 14).Value = -5893.5
$ws.Cells.Item(100, 8).Value = 9525545
$ws.Cells.Item(100, 9).Value = 12346775
$ws.Cells.Item(100, 10).Value = 3894.625
$ws.Cells.Item(100, 11).Value = 12346775
$ws.Cells.Item(100, 12).Value = 3894.625
$ws.Cells.Item(100, 13).Value = -12346234
$ws.Cells.Item(100, 14).Value = -4976.625
$ws.Cells.Item(135, 8).Value = 3345.28
$ws.Cells.Item(135, 9).Value = 2108.4443
$ws.Cells.Item(135, 10).Value = 6525.7144
$ws.Cells.Item(135, 11).Value = 18975.9987
$ws.Cells.Item(135, 12).Value = 58731.4296
$ws.Cells.Item(135, 13).Value = -16440.9987
$ws.Cells.Item(135, 14).Value = -63801.4296
$ws.Cells.Item(137, 8).Value = 1758.5
$ws.Cells.Item(137, 9).Value = 1487.75
$ws.Cells.Item(137, 10).Value = 2300
$ws.Cells.Item(137, 11).Value = 4463.25
$ws.Cells.Item(137, 12).Value = 6900
$ws.Cells.Item(137, 13).Value = -1913.25
$ws.Cells.Item(137, 14).Value = -12000

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 130
$ws.Cells.Item(5, 9).Value = 137.5
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 137.5
$ws.Cells.Item(5, 12).Value = 100
$ws.Cells.Item(5, 13).Value = -25.5
$ws.Cells.Item(5, 14).Value = -324
$ws.Cells.Item(61, 8).Value = 2943.0625
$ws.Cells.Item(61, 9).Value = 3086.8845
$ws.Cells.Item(61, 11).Value = 3086.8845
$ws.Cells.Item(61, 13).Value = -2874.8845
$ws.Cells.Item(74, 8).Value = 7990
$ws.Cells.Item(74, 9).Value = 15475
$ws.Cells.Item(74, 11).Value = 15475
$ws.Cells.Item(74, 13).Value = -14601
$ws.Cells.Item(77, 8).Value = 7990
$ws.Cells.Item(77, 9).Value = 15475
$ws.Cells.Item(77, 11).Value = 77375
$ws.Cells.Item(77, 13).Value = -73007
$ws.Cells.Item(97, 8).Value = 2104.2144
$ws.Cells.Item(97, 9).Value = 2001.6111
$ws.Cells.Item(97, 10).Value = 2288.9
$ws.Cells.Item(97, 11).Value = 2001.6111
$ws.Cells.Item(97, 12).Value = 2288.9
$ws.Cells.Item(97, 13).Value = -1505.6111
$ws.Cells.Item(97, 14).Value = -3280.9
$ws.Cells.Item(136, 8).Value = 2943.0625
$ws.Cells.Item(136, 9).Value = 3086.8845
$ws.Cells.Item(136, 11).Value = 9260.6535
$ws.Cells.Item(136, 13).Value = -6710.6535

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 130
$ws.Cells.Item(4, 9).Value = 137.5
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 137.5
$ws.Cells.Item(4, 12).Value = 100
$ws.Cells.Item(4, 13).Value = -22.5
$ws.Cells.Item(4, 14).Value = -330
$ws.Cells.Item(86, 8).Value = 47624416
$ws.Cells.Item(86, 9).Value = 100001940
$ws.Cells.Item(86, 10).Value = 8482.454
$ws.Cells.Item(86, 11).Value = 100001940
$ws.Cells.Item(86, 12).Value = 8482.454
$ws.Cells.Item(86, 13).Value = -100000817
$ws.Cells.Item(86, 14).Value = -10728.454
$ws.Cells.Item(89, 8).Value = 47624416
$ws.Cells.Item(89, 9).Value = 100001940
$ws.Cells.Item(89, 10).Value = 8482.454
$ws.Cells.Item(89, 11).Value = 500009700
$ws.Cells.Item(89, 12).Value = 42412.27
$ws.Cells.Item(89, 13).Value = -500004084
$ws.Cells.Item(89, 14).Value = -53644.27

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4633.0527
$ws.Cells.Item(31, 9).Value = 3272.4285
$ws.Cells.Item(31, 10).Value = 8442.8
$ws.Cells.Item(31, 11).Value = 3272.4285
$ws.Cells.Item(31, 12).Value = 8442.8
$ws.Cells.Item(31, 13).Value = -2977.4285
$ws.Cells.Item(31, 14).Value = -9032.8
$ws.Cells.Item(34, 8).Value = 4633.0527
$ws.Cells.Item(34, 9).Value = 3272.4285
$ws.Cells.Item(34, 10).Value = 8442.8
$ws.Cells.Item(34, 11).Value = 3272.4285
$ws.Cells.Item(34, 12).Value = 8442.8
$ws.Cells.Item(34, 13).Value = -3070.4285
$ws.Cells.Item(34, 14).Value = -8846.8
$ws.Cells.Item(94, 8).Value = 5772.857
$ws.Cells.Item(94, 10).Value = 5772.857
$ws.Cells.Item(94, 12).Value = 5772.857
$ws.Cells.Item(94, 14).Value = -6674.857
$ws.Cells.Item(95, 8).Value = 30540
$ws.Cells.Item(95, 10).Value = 30540
$ws.Cells.Item(95, 12).Value = 30540
$ws.Cells.Item(95, 14).Value = -36032
$ws.Cells.Item(105, 8).Value = 1723
$ws.Cells.Item(105, 9).Value = 1318.5714
$ws.Cells.Item(105, 11).Value = 1318.5714
$ws.Cells.Item(105, 13).Value = 428.4286
$ws.Cells.Item(108, 8).Value = 35813
$ws.Cells.Item(108, 10).Value = 35813
$ws.Cells.Item(108, 12).Value = 35813
$ws.Cells.Item(108, 14).Value = -43493

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1187.375
$ws.Cells.Item(68, 9).Value = 1000
$ws.Cells.Item(68, 10).Value = 1374.75
$ws.Cells.Item(68, 11).Value = 3000
$ws.Cells.Item(68, 12).Value = 4124.25
$ws.Cells.Item(68, 13).Value = -2189
$ws.Cells.Item(68, 14).Value = -5746.25
$ws.Cells.Item(71, 8).Value = 1187.375
$ws.Cells.Item(71, 9).Value = 1000
$ws.Cells.Item(71, 10).Value = 1374.75
$ws.Cells.Item(71, 11).Value = 9000
$ws.Cells.Item(71, 12).Value = 12372.75
$ws.Cells.Item(71, 13).Value = -4944
$ws.Cells.Item(71, 14).Value = -20484.75
$ws.Cells.Item(122, 8).Value = 739.3043
$ws.Cells.Item(122, 9).Value = 290.15384
$ws.Cells.Item(122, 10).Value = 1323.2
$ws.Cells.Item(122, 11).Value = 2611.38456
$ws.Cells.Item(122, 12).Value = 11908.8
$ws.Cells.Item(122, 13).Value = -161.38456
$ws.Cells.Item(122, 14).Value = -16808.8
$ws.Cells.Item(131, 8).Value = 826.17
$ws.Cells.Item(131, 10).Value = 848.81055
$ws.Cells.Item(131, 12).Value = 2546.43165
$ws.Cells.Item(131, 14).Value = -12626.43165

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2204.2856
$ws.Cells.Item(102, 9).Value = 2129.5881
$ws.Cells.Item(102, 10).Value = 2521.75
$ws.Cells.Item(102, 11).Value = 2129.5881
$ws.Cells.Item(102, 12).Value = 2521.75
$ws.Cells.Item(102, 13).Value = -507.5880999999999
$ws.Cells.Item(102, 14).Value = -5765.75
$ws.Cells.Item(113, 8).Value = 66668564
$ws.Cells.Item(113, 9).Value = 1850.625
$ws.Cells.Item(113, 10).Value = 142859090
$ws.Cells.Item(113, 11).Value = 1850.625
$ws.Cells.Item(113, 12).Value = 142859090
$ws.Cells.Item(113, 13).Value = 319.375
$ws.Cells.Item(113, 14).Value = -142863430
$ws.Cells.Item(122, 8).Value = 2004.1904
$ws.Cells.Item(122, 9).Value = 1911.1765
$ws.Cells.Item(122, 10).Value = 2399.5
$ws.Cells.Item(122, 11).Value = 5733.529500000001
$ws.Cells.Item(122, 12).Value = 7198.5
$ws.Cells.Item(122, 13).Value = -3283.529500000001
$ws.Cells.Item(122, 14).Value = -12098.5
$ws.Cells.Item(126, 8).Value = 15154405
$ws.Cells.Item(126, 9).Value = 3331.111
$ws.Cells.Item(126, 10).Value = 83334240
$ws.Cells.Item(126, 11).Value = 9993.332999999999
$ws.Cells.Item(126, 12).Value = 250002720
$ws.Cells.Item(126, 13).Value = -7523.332999999999
$ws.Cells.Item(126, 14).Value = -250007660
$ws.Cells.Item(132, 8).Value = 2124.6572
$ws.Cells.Item(132, 9).Value = 1527.125
$ws.Cells.Item(132, 10).Value = 3428.3635
$ws.Cells.Item(132, 11).Value = 4581.375
$ws.Cells.Item(132, 12).Value = 10285.0905
$ws.Cells.Item(132, 13).Value = -2051.375
$ws.Cells.Item(132, 14).Value = -15345.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1841.079
$ws.Cells.Item(40, 9).Value = 1591.2084
$ws.Cells.Item(40, 10).Value = 2269.4285
$ws.Cells.Item(40, 11).Value = 1591.2084
$ws.Cells.Item(40, 12).Value = 2269.4285
$ws.Cells.Item(40, 13).Value = -1455.2084
$ws.Cells.Item(40, 14).Value = -2541.4285
$ws.Cells.Item(46, 8).Value = 11680.1
$ws.Cells.Item(46, 9).Value = 2425.25
$ws.Cells.Item(46, 10).Value = 17850
$ws.Cells.Item(46, 11).Value = 2425.25
$ws.Cells.Item(46, 12).Value = 17850
$ws.Cells.Item(46, 13).Value = -2237.25
$ws.Cells.Item(46, 14).Value = -18226
$ws.Cells.Item(68, 8).Value = 1816.4054
$ws.Cells.Item(68, 9).Value = 1771
$ws.Cells.Item(68, 10).Value = 1900.2307
$ws.Cells.Item(68, 11).Value = 1771
$ws.Cells.Item(68, 12).Value = 1900.2307
$ws.Cells.Item(68, 13).Value = -1022
$ws.Cells.Item(68, 14).Value = -3398.2307
$ws.Cells.Item(71, 8).Value = 1816.4054
$ws.Cells.Item(71, 9).Value = 1771
$ws.Cells.Item(71, 10).Value = 1900.2307
$ws.Cells.Item(71, 11).Value = 8855
$ws.Cells.Item(71, 12).Value = 9501.1535
$ws.Cells.Item(71, 13).Value = -5111
$ws.Cells.Item(71, 14).Value = -16989.1535
$ws.Cells.Item(136, 8).Value = 4028.8572
$ws.Cells.Item(136, 9).Value = 2860.5
$ws.Cells.Item(136, 10).Value = 5586.6665
$ws.Cells.Item(136, 11).Value = 8581.5
$ws.Cells.Item(136, 12).Value = 16759.9995
$ws.Cells.Item(136, 13).Value = -6031.5
$ws.Cells.Item(136, 14).Value = -21859.9995
$ws.Cells.Item(141, 8).Value = 67987.5
$ws.Cells.Item(141, 10).Value = 67987.5
$ws.Cells.Item(141, 12).Value = 67987.5
$ws.Cells.Item(141, 14).Value = -78347.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 815.6
$ws.Cells.Item(107, 9).Value = 784
$ws.Cells.Item(107, 10).Value = 1100
$ws.Cells.Item(107, 11).Value = 2352
$ws.Cells.Item(107, 12).Value = 3300
$ws.Cells.Item(107, 13).Value = -432
$ws.Cells.Item(107, 14).Value = -7140
$ws.Cells.Item(132, 8).Value = 34887570
$ws.Cells.Item(132, 9).Value = 51725530
$ws.Cells.Item(132, 10).Value = 8938.429
$ws.Cells.Item(132, 11).Value = 155176590
$ws.Cells.Item(132, 12).Value = 26815.287
$ws.Cells.Item(132, 13).Value = -155174060
$ws.Cells.Item(132, 14).Value = -31875.287
$ws.Cells.Item(136, 8).Value = 1130.0968
$ws.Cells.Item(136, 9).Value = 905.1667
$ws.Cells.Item(136, 10).Value = 1441.5385
$ws.Cells.Item(136, 11).Value = 2715.5001
$ws.Cells.Item(136, 12).Value = 4324.6155
$ws.Cells.Item(136, 13).Value = -165.5001000000002
$ws.Cells.Item(136, 14).Value = -9424.6155
